$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Rewrite the data table ----
# Header row (unchanged text, just shared-string index churn upstream)
$ws.Range("A1").Value = "function"
$ws.Range("B1").Value = "area"
$ws.Range("C1").Value = "area_outdoors"
$ws.Range("D1").Value = "daylight_analysis"
$ws.Range("E1").Value = "sunlight_analysis"
$ws.Range("F1").Value = "height_analysis"
$ws.Range("G1").Value = "inv_height_analysis"
$ws.Range("H1").Value = "noise"

$ws.Range("A2").Value = "str"
$ws.Range("B2").Value = "int"
$ws.Range("C2").Value = "int"
$ws.Range("D2").Value = "float"
$ws.Range("E2").Value = "float"
$ws.Range("F2").Value = "float"
$ws.Range("G2").Value = "float"
$ws.Range("H2").Value = "float"

# row 3: student_housing
$ws.Range("A3").Value = "student_housing"
$ws.Range("B3").Value = 10000
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = "0.5"
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = "0.5"

# row 4: starter_housing
$ws.Range("A4").Value = "starter_housing"
$ws.Range("B4").Value = 12000
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = "0.5"
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = "0.7"

# row 5: elderly_housing
$ws.Range("A5").Value = "elderly_housing"
$ws.Range("B5").Value = 13000
$ws.Range("C5").Value = 1000
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = "0.5"
$ws.Range("F5").Value = "0.5"
$ws.Range("G5").Value = "0.5"
$ws.Range("H5").Value = "0.9"

# row 6: sportfacilities (was carparking_private)
$ws.Range("A6").Value = "sportfacilities"
$ws.Range("B6").Value = 600
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = "0.6"
$ws.Range("E6").Value = "0.3"
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 1
$ws.Range("H6").Value = "0.3"

# row 7: kindergarden (was bikeparking_private)
$ws.Range("A7").Value = "kindergarden"
$ws.Range("B7").Value = 200
$ws.Range("C7").Value = 240
$ws.Range("D7").Value = "0.8"
$ws.Range("E7").Value = "0.7"
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 1
$ws.Range("H7").Value = "0.8"

# row 8: care_centre (was carparking_public formula row)
$ws.Range("A8").Value = "care_centre"
$ws.Range("B8").ClearContents()
$ws.Range("B8").Value = 230
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = "0.8"
$ws.Range("E8").Value = "0.3"
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 1
$ws.Range("H8").Value = "0.8"

# row 9: physical_therapy_centre (was bikeparking_public formula row)
$ws.Range("A9").Value = "physical_therapy_centre"
$ws.Range("B9").ClearContents()
$ws.Range("B9").Value = 60
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = "0.8"
$ws.Range("E9").Value = "0.3"
$ws.Range("F9").Value = 0
$ws.Range("G9").Value = 1
$ws.Range("H9").Value = "0.8"

# row 10: working_places (was sportfacilities)
$ws.Range("A10").Value = "working_places"
$ws.Range("B10").Value = 550
$ws.Range("C10").Value = 0
$ws.Range("D10").Value = "0.7"
$ws.Range("E10").Value = "0.3"
$ws.Range("F10").Value = 0
$ws.Range("G10").Value = 1
$ws.Range("H10").Value = 1

# row 11: artstudios (was kindergarden)
$ws.Range("A11").Value = "artstudios"
$ws.Range("B11").Value = 60
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = "0.8"
$ws.Range("E11").Value = "0.6"
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = 1
$ws.Range("H11").Value = "0.8"

# row 12: library_cafe (was city_garden)
$ws.Range("A12").Value = "library_cafe"
$ws.Range("B12").Value = 900
$ws.Range("C12").Value = 0
$ws.Range("D12").Value = "0.8"
$ws.Range("E12").Value = "0.6"
$ws.Range("F12").Value = 0
$ws.Range("G12").Value = 1
$ws.Range("H12").Value = "0.9"

# row 13: community_spaces_student_starters (was care_centre)
$ws.Range("A13").Value = "community_spaces_student_starters"
$ws.Range("B13").Value = 500
$ws.Range("C13").Value = 0
$ws.Range("D13").Value = "0.8"
$ws.Range("E13").Value = "0.3"
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = "0.7"
$ws.Range("H13").Value = "0.7"

# row 14: community_spaces_elderly (was physical_therapy_centre)
$ws.Range("A14").Value = "community_spaces_elderly"
$ws.Range("B14").Value = 300
$ws.Range("C14").Value = 0
$ws.Range("D14").Value = "0.8"
$ws.Range("E14").Value = "0.3"
$ws.Range("F14").Value = 0
$ws.Range("G14").Value = "0.7"
$ws.Range("H14").Value = "0.7"

# row 15: laundry_room (brand-new row, was working_places)
$ws.Range("A15").Value = "laundry_room"
$ws.Range("B15").Value = 100
$ws.Range("C15").Value = 0
$ws.Range("D15").Value = "0.5"
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 0

# Old rows 16-19 (artstudios, library_cafe, community_spaces_student_starters,
# community_spaces_elderly) are gone as named rows; clear their previous content.
$ws.Range("A16:A19").ClearContents()

# rows 16-37: blank placeholder rows (formatted, no values) in columns B:H
for ($r = 16; $r -le 37; $r++) {
    $ws.Range("B$r:H$r").ClearContents()
}

$ws.Range("I10").Select()
